# Update cryptos list with refreshed Price (D) / Volume(1h) (E) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '74.961.60'
$ws.Range('E2').Value = '  +1.52%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.818.18'
$ws.Range('E3').Value = '  +7.31%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.00%  '

# Row 5: Solana
$ws.Range('D5').Value = '''187.33'
$ws.Range('E5').Value = '  +1.19%  '

# Row 6: BNB
$ws.Range('D6').Value = '''594.56'

# Row 7: USDC
$ws.Range('E7').Value = '  +0.00%  '

# Row 8: XRP
$ws.Range('D8').Value = '''0.550'
$ws.Range('E8').Value = '  +2.87%  '

# Row 9: Dogecoin
$ws.Range('D9').Value = '''0.192'
$ws.Range('E9').Value = '  -4.65%  '

# Row 10: LidoStakedEther
$ws.Range('D10').Value = '2.817.79'
$ws.Range('E10').Value = '  +7.39%  '

# Row 11: TRON
$ws.Range('E11').Value = '  -1.24%  '

# Row 12: Cardano
$ws.Range('E12').Value = '  +3.25%  '

# Row 13: Toncoin
$ws.Range('D13').Value = '''4.87'
$ws.Range('E13').Value = '  +2.12%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '3.337.55'
$ws.Range('E14').Value = '  +7.36%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '74.895.01'
$ws.Range('E15').Value = '  +1.60%  '

# Row 16: ShibaInu
$ws.Range('E16').Value = '  -0.78%  '

# Row 17: Avalanche
$ws.Range('D17').Value = '''26.79'
$ws.Range('E17').Value = '  +2.35%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '2.818.64'
$ws.Range('E18').Value = '  +6.80%  '

# Row 19: Uniswap
$ws.Range('E19').Value = '  -1.07%  '

# Row 20: Chainlink
$ws.Range('D20').Value = '''12.28'
$ws.Range('E20').Value = '  +3.84%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''376.54'
$ws.Range('E21').Value = '  +1.34%  '

# Row 22: SuiNetwork
$ws.Range('D22').Value = '''2.24'
$ws.Range('E22').Value = '  -2.29%  '

# Row 23: Polkadot
$ws.Range('E23').Value = '  -0.80%  '

# Row 24: LEO
$ws.Range('E24').Value = '  -0.28%  '

# Row 25: Dai
$ws.Range('D25').Value = '''0.999'
$ws.Range('E25').Value = '  -0.14%  '

# Row 26: Litecoin
$ws.Range('D26').Value = '''70.67'
$ws.Range('E26').Value = '  +1.20%  '

# Row 27: WrappedeETH
$ws.Range('D27').Value = '2.961.02'
$ws.Range('E27').Value = '  +7.19%  '

# Row 28: NEARProtocol
$ws.Range('D28').Value = '''4.16'
$ws.Range('E28').Value = '  +0.34%  '

# Row 29: Aptos
$ws.Range('D29').Value = '''9.68'
$ws.Range('E29').Value = '  +3.17%  '

# Row 30: PEPE
$ws.Range('E30').Value = '  +10.27%  '

# Row 31: Binance-PegBSC-USD
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  -0.52%  '

# Row 32: Bittensor
$ws.Range('D32').Value = '''514.77'
$ws.Range('E32').Value = '  -1.57%  '

# Row 33: Fetch.AI
$ws.Range('E33').Value = '  +0.01%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').Value = '''7.71'
$ws.Range('E34').Value = '  +0.73%  '

# Row 35: PancakeSwap
$ws.Range('E35').Value = '  +2.40%  '

# Row 37: Monero
$ws.Range('D37').Value = '''162.40'
$ws.Range('E37').Value = '  +0.21%  '

# Row 38: EthereumClassic
$ws.Range('E38').Value = '  +3.94%  '

# Row 39: Kaspa
$ws.Range('E39').Value = '  -1.17%  '

# Row 40: WhiteBITCoin
$ws.Range('E40').Value = '  +0.63%  '

# Row 41: Aave
$ws.Range('D41').Value = '''185.59'
$ws.Range('E41').Value = '  +15.07%  '

# Row 42: USDe
$ws.Range('E42').Value = '  +0.01%  '

# Row 43: PolygonEcosystemToken
$ws.Range('D43').Value = '''0.339'
$ws.Range('E43').Value = '  +3.03%  '

# Row 44: RenderToken
$ws.Range('D44').Value = '''4.99'
$ws.Range('E44').Value = '  +1.59%  '

# Row 45: Stacks
$ws.Range('E45').Value = '  -0.38%  '

# Row 46: ImmutableX
$ws.Range('E46').Value = '  +1.50%  '

# Row 47: OKB
$ws.Range('D47').Value = '''39.97'
$ws.Range('E47').Value = '  +2.60%  '

# Row 48: dogwifhat
$ws.Range('D48').Value = '''2.33'
$ws.Range('E48').Value = '  -1.55%  '

# Row 49: Cronos
$ws.Range('E49').Value = '  -0.03%  '

# Row 50: ARBITRUM
$ws.Range('D50').Value = '''0.571'
$ws.Range('E50').Value = '  +7.85%  '

# Row 51: Filecoin
$ws.Range('D51').Value = '''3.71'
$ws.Range('E51').Value = '  +2.21%  '
